$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure price column (D) cells are stored as text, matching the
# original workbook's inline-string cells (avoids "1.00" -> 1 numeric coercion).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Cells.Item(2, 4).Value = "59.386.13"
$ws.Cells.Item(2, 5).Value = "  -5.62%  "
$ws.Cells.Item(3, 4).Value = "2.458.67"
$ws.Cells.Item(3, 5).Value = "  -8.21%  "
$ws.Cells.Item(4, 4).Value = "1.00"
$ws.Cells.Item(4, 5).Value = "  -0.02%  "
$ws.Cells.Item(5, 4).Value = "539.47"
$ws.Cells.Item(5, 5).Value = "  -2.48%  "
$ws.Cells.Item(6, 4).Value = "148.07"
$ws.Cells.Item(6, 5).Value = "  -6.41%  "
$ws.Cells.Item(7, 5).Value = "  -0.26%  "
$ws.Cells.Item(8, 4).Value = "0.569"
$ws.Cells.Item(8, 5).Value = "  -4.10%  "
$ws.Cells.Item(9, 4).Value = "2.477.83"
$ws.Cells.Item(9, 5).Value = "  -7.60%  "
$ws.Cells.Item(10, 4).Value = "0.0995"
$ws.Cells.Item(10, 5).Value = "  -5.74%  "
$ws.Cells.Item(11, 5).Value = "  -2.59%  "
$ws.Cells.Item(12, 5).Value = "  -0.63%  "
$ws.Cells.Item(13, 4).Value = "0.353"
$ws.Cells.Item(13, 5).Value = "  -3.95%  "
$ws.Cells.Item(14, 4).Value = "2.900.17"
$ws.Cells.Item(14, 5).Value = "  -8.03%  "
$ws.Cells.Item(15, 4).Value = "24.18"
$ws.Cells.Item(15, 5).Value = "  -7.69%  "
$ws.Cells.Item(16, 4).Value = "59.332.25"
$ws.Cells.Item(17, 5).Value = "  -5.58%  "
$ws.Cells.Item(18, 4).Value = "2.519.71"
$ws.Cells.Item(18, 5).Value = "  -6.04%  "
$ws.Cells.Item(19, 4).Value = "11.18"
$ws.Cells.Item(19, 5).Value = "  -5.70%  "
$ws.Cells.Item(20, 4).Value = "4.35"
$ws.Cells.Item(20, 5).Value = "  -5.47%  "
$ws.Cells.Item(21, 4).Value = "324.68"
$ws.Cells.Item(21, 5).Value = "  -5.82%  "
$ws.Cells.Item(22, 4).Value = "0.969"
$ws.Cells.Item(22, 5).Value = "  -2.95%  "
$ws.Cells.Item(23, 4).Value = "5.75"
$ws.Cells.Item(23, 5).Value = "  -8.61%  "
$ws.Cells.Item(24, 4).Value = "0.462"
$ws.Cells.Item(24, 5).Value = "  -8.85%  "
$ws.Cells.Item(25, 4).Value = "60.60"
$ws.Cells.Item(25, 5).Value = "  -4.36%  "
$ws.Cells.Item(26, 5).Value = "  -3.65%  "
$ws.Cells.Item(27, 5).Value = "  -2.06%  "
$ws.Cells.Item(28, 4).Value = "7.72"
$ws.Cells.Item(28, 5).Value = "  -5.63%  "
$ws.Cells.Item(29, 4).Value = "6.77"
$ws.Cells.Item(29, 5).Value = "  -6.47%  "
$ws.Cells.Item(30, 2).Value = "Fetch.AI"
$ws.Cells.Item(30, 3).Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Cells.Item(30, 4).Value = "1.27"
$ws.Cells.Item(30, 5).Value = "  -7.94%  "
$ws.Cells.Item(31, 2).Value = "PancakeSwap"
$ws.Cells.Item(31, 3).Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Cells.Item(31, 4).Value = "1.83"
$ws.Cells.Item(31, 5).Value = "  -6.07%  "
$ws.Cells.Item(32, 4).Value = "0.0₃0775"
$ws.Cells.Item(32, 5).Value = "  -9.47%  "
$ws.Cells.Item(33, 5).Value = "  -0.14%  "
$ws.Cells.Item(34, 4).Value = "158.47"
$ws.Cells.Item(34, 5).Value = "  -4.04%  "
$ws.Cells.Item(35, 4).Value = "4.60"
$ws.Cells.Item(35, 5).Value = "  -5.41%  "
$ws.Cells.Item(36, 5).Value = "  -6.89%  "
$ws.Cells.Item(37, 4).Value = "18.40"
$ws.Cells.Item(37, 5).Value = "  -5.60%  "
$ws.Cells.Item(38, 4).Value = "1.77"
$ws.Cells.Item(38, 5).Value = "  -0.66%  "
$ws.Cells.Item(39, 4).Value = "5.98"
$ws.Cells.Item(39, 5).Value = "  -5.57%  "
$ws.Cells.Item(40, 4).Value = "322.45"
$ws.Cells.Item(40, 5).Value = "  -7.69%  "
$ws.Cells.Item(41, 4).Value = "36.82"
$ws.Cells.Item(41, 5).Value = "  -3.92%  "
$ws.Cells.Item(42, 4).Value = "0.837"
$ws.Cells.Item(42, 5).Value = "  -12.89%  "
$ws.Cells.Item(43, 5).Value = "  -7.26%  "
$ws.Cells.Item(44, 4).Value = "0.995"
$ws.Cells.Item(44, 5).Value = "  -0.36%  "
$ws.Cells.Item(45, 4).Value = "10.71"
$ws.Cells.Item(45, 5).Value = "  -2.86%  "
$ws.Cells.Item(46, 5).Value = "  -4.87%  "
$ws.Cells.Item(47, 4).Value = "0.0941"
$ws.Cells.Item(47, 5).Value = "  -3.14%  "
$ws.Cells.Item(48, 4).Value = "0.0526"
$ws.Cells.Item(48, 5).Value = "  -6.14%  "
$ws.Cells.Item(49, 4).Value = "19.08"
$ws.Cells.Item(49, 5).Value = "  -8.21%  "
$ws.Cells.Item(50, 4).Value = "18.58"
$ws.Cells.Item(50, 5).Value = "  -8.84%  "
$ws.Cells.Item(51, 4).Value = "122.22"
$ws.Cells.Item(51, 5).Value = "  -4.88%  "
